$wb = $excel.ActiveWorkbook

# The new sheet is a continuation of the existing position-statistics sheets.
# Duplicate the last sheet ("20191027") to the end of the workbook and
# rename it to "20191029", then update it with the new day's numbers.
$srcSheet = $wb.Worksheets.Item("20191027")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "20191029"

# Update the "B" (count) column for rows 1-14, and a couple of "C" values.
$ws.Range("B1").Value = 877
$ws.Range("B2").Value = 382
$ws.Range("B3").Value = 453
$ws.Range("B4").Value = 512
$ws.Range("B5").Value = 567
$ws.Range("B6").Value = 599
$ws.Range("B7").Value = 892
$ws.Range("B8").Value = 774
$ws.Range("B9").Value = 872
$ws.Range("B10").Value = 1125
$ws.Range("B11").Value = 1802
$ws.Range("B12").Value = 1759
$ws.Range("B13").Value = 1149
$ws.Range("C13").Value = 0.08
$ws.Range("B14").Value = 2344
$ws.Range("C14").Value = 0.16

# Update the bottom "look multi/short/flat" block.
$ws.Range("B20").Value = 3755
$ws.Range("A21").Value = "看空 (已选)"
$ws.Range("B21").Value = 6288
$ws.Range("A22").Value = "看平"
$ws.Range("B22").Value = 1761
$ws.Range("B23").Value = 2325

# Move the selection/active cell on the new sheet.
$ws.Range("M17").Select()

# Make the new sheet the active tab.
$ws.Activate()
